# Generate Report for Handoff
# The handoff-report generation process re-ran and produced a fresh
# "Latest Handoff Datetime" for the 86c2c410-695a-4bac-adb4-97f783701b8c
# source file entry on the zh-cn status sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("H4").Value = "2017-02-21 03:41:25"
